$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Move sensitivities to targets": the overprovision-of-food sensitivity
# parameters (overFmean / overFdelta) are expressed per day, not as a
# percentage, so update their Unit column accordingly.
$ws.Range("D8").Value = "d-1"
$ws.Range("D9").Value = "d-1"

# Move the active selection to C11 (the empty target cell below the table)
$ws.Range("C11").Select()
